# Auto-generated script applying market-data refresh to Kraken_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 91.666664
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 92.5
$ws.Range("K2").Value = 90
$ws.Range("L2").Value = 92.5
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = -318.5

$ws.Range("H28").Value = 2042.1904
$ws.Range("I28").Value = 2066.8235
$ws.Range("J28").Value = 1937.5
$ws.Range("K28").Value = 2066.8235
$ws.Range("L28").Value = 1937.5
$ws.Range("M28").Value = -1581.8235
$ws.Range("N28").Value = -2907.5

$ws.Range("H41").Value = 2765.3
$ws.Range("I41").Value = 199
$ws.Range("J41").Value = 3406.875
$ws.Range("K41").Value = 199
$ws.Range("L41").Value = 3406.875
$ws.Range("M41").Value = 241
$ws.Range("N41").Value = -4286.875

$ws.Range("H53").Value = 113.375
$ws.Range("I53").Value = 115.28571
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 115.28571
$ws.Range("L53").Value = 100
$ws.Range("M53").Value = 521.71429
$ws.Range("N53").Value = -1374

$ws.Range("H62").Value = 10750
$ws.Range("I62").Value = 10750
$ws.Range("K62").Value = 10750
$ws.Range("M62").Value = -10126

$ws.Range("H64").Value = 4665
$ws.Range("J64").Value = 4665
$ws.Range("L64").Value = 4665
$ws.Range("N64").Value = -5161

$ws.Range("H65").Value = 10750
$ws.Range("I65").Value = 10750
$ws.Range("K65").Value = 53750
$ws.Range("M65").Value = -50630

$ws.Range("H67").Value = 4665
$ws.Range("J67").Value = 4665
$ws.Range("L67").Value = 4665
$ws.Range("N67").Value = -6381

$ws.Range("H70").Value = 1823.75
$ws.Range("I70").Value = 1330
$ws.Range("J70").Value = 2120
$ws.Range("K70").Value = 3990
$ws.Range("L70").Value = 6360
$ws.Range("M70").Value = -3720
$ws.Range("N70").Value = -6900

$ws.Range("H73").Value = 1823.75
$ws.Range("I73").Value = 1330
$ws.Range("J73").Value = 2120
$ws.Range("K73").Value = 3990
$ws.Range("L73").Value = 6360
$ws.Range("M73").Value = -3054
$ws.Range("N73").Value = -8232

$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("K76").Value = 2000
$ws.Range("M76").Value = -1685

$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("K79").Value = 2000
$ws.Range("M79").Value = -908

$ws.Range("H86").Value = 8500
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877

$ws.Range("H89").Value = 8500
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384

$ws.Range("H92").Value = 424.57144
$ws.Range("I92").Value = 210.75
$ws.Range("K92").Value = 210.75
$ws.Range("M92").Value = 1037.25

$ws.Range("H98").Value = 4692.5
$ws.Range("I98").Value = 4692.5
$ws.Range("K98").Value = 4692.5
$ws.Range("M98").Value = -3194.5

$ws.Range("H106").Value = 1827.5714
$ws.Range("I106").Value = 1298.8334
$ws.Range("K106").Value = 1298.8334
$ws.Range("M106").Value = -667.8334

$ws.Range("H107").Value = 1449.2222
$ws.Range("I107").Value = 649.1429000000001
$ws.Range("J107").Value = 4249.5
$ws.Range("K107").Value = 649.1429000000001
$ws.Range("L107").Value = 4249.5
$ws.Range("M107").Value = 1270.8571
$ws.Range("N107").Value = -8089.5

$ws.Range("H122").Value = 4692.5
$ws.Range("I122").Value = 4692.5
$ws.Range("K122").Value = 14077.5
$ws.Range("M122").Value = -11627.5

$ws.Range("H132").Value = 4129.5625
$ws.Range("I132").Value = 1688.5454
$ws.Range("K132").Value = 5065.6362
$ws.Range("M132").Value = -2535.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1665.3334
$ws.Range("I25").Value = 1665.3334
$ws.Range("K25").Value = 1665.3334
$ws.Range("M25").Value = -1263.3334

$ws.Range("H63").Value = 2598.8
$ws.Range("J63").Value = 1998.6666
$ws.Range("L63").Value = 1998.6666
$ws.Range("N63").Value = -3370.6666

$ws.Range("H66").Value = 2598.8
$ws.Range("J66").Value = 1998.6666
$ws.Range("L66").Value = 9993.333000000001
$ws.Range("N66").Value = -16857.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 7500
$ws.Range("I11").Value = 10000
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -9860
$ws.Range("N11").Value = -5280

$ws.Range("H37").Value = 349
$ws.Range("I37").Value = 349
$ws.Range("K37").Value = 349
$ws.Range("M37").Value = -212

$ws.Range("H82").Value = 21362.5
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766

$ws.Range("H85").Value = 21362.5
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652

$ws.Range("H102").Value = 34999
$ws.Range("I102").Value = 34999
$ws.Range("K102").Value = 34999
$ws.Range("M102").Value = -31754

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6311.375
$ws.Range("I31").Value = 5960.8887
$ws.Range("K31").Value = 5960.8887
$ws.Range("M31").Value = -5665.8887

$ws.Range("H34").Value = 6311.375
$ws.Range("I34").Value = 5960.8887
$ws.Range("K34").Value = 5960.8887
$ws.Range("M34").Value = -5758.8887

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3000
$ws.Range("I3").Value = 3000
$ws.Range("K3").Value = 9000
$ws.Range("M3").Value = -8888

$ws.Range("H121").Value = 1611.2222
$ws.Range("I121").Value = 359.33334
$ws.Range("J121").Value = 2237.1667
$ws.Range("K121").Value = 1078.00002
$ws.Range("L121").Value = 6711.500100000001
$ws.Range("M121").Value = 231.9999800000001
$ws.Range("N121").Value = -9331.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3475.6
$ws.Range("I102").Value = 3475.6
$ws.Range("K102").Value = 3475.6
$ws.Range("M102").Value = -1853.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3028
$ws.Range("I68").Value = 2999.3333
$ws.Range("K68").Value = 2999.3333
$ws.Range("M68").Value = -2250.3333

$ws.Range("H71").Value = 3028
$ws.Range("I71").Value = 2999.3333
$ws.Range("K71").Value = 14996.6665
$ws.Range("M71").Value = -11252.6665

